# Convert M2Doc field-code runs (fldChar begin / instrText / fldChar end)
# into plain-text "{...}" runs, e.g.
#   { begin } { instrText:  m:Sequence{...}  } { end }
# becomes
#   { t: "{m:Sequence{...}}" }
#
# Mirrors the parser change to TokenIteratorFieldRewriterSplit: M2Doc
# tokens no longer need to live inside a real Word field, a simple
# "{...}" text run is enough.

$d = $word.ActiveDocument

# Walk fields back-to-front so deleting one does not shift the
# index/position of the fields still to be processed.
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $fld = $d.Fields.Item($i)

    # Grab the field's instruction text (the M2Doc token) before the
    # field itself is removed.
    $codeText = $fld.Code.Text
    $codeStart = $fld.Code.Start

    # Locate the start of the paragraph that hosts the field so the
    # replacement text can be scoped to just that paragraph (leaving
    # sibling paragraphs untouched).
    $paraStart = -1
    for ($p = 1; $p -le $d.Paragraphs.Count; $p++) {
        $pr = $d.Paragraphs.Item($p).Range
        if ($codeStart -ge $pr.Start -and $codeStart -lt $pr.End) {
            $paraStart = $pr.Start
            break
        }
    }

    # Trim the leading/trailing space Word stores around field codes and
    # wrap the token in the "{ ... }" textual-token syntax.
    $newText = "{" + $codeText.Trim() + "}"

    # Remove the begin/instrText/end run triplet entirely. This shifts
    # every position at/after the field, so the paragraph's end (and any
    # position captured before the delete) can no longer be trusted.
    $fld.Delete()

    # Re-locate the now-empty paragraph via its (still valid) start
    # offset and grab its fresh, post-delete end offset.
    $paraEnd = -1
    for ($p = 1; $p -le $d.Paragraphs.Count; $p++) {
        $pr = $d.Paragraphs.Item($p).Range
        if ($pr.Start -eq $paraStart) {
            $paraEnd = $pr.End
            break
        }
    }

    # Replace what is left (an empty run) with the plain-text token.
    # Exclude the trailing paragraph-mark position (paraEnd - 1) so the
    # new text lands inside the paragraph instead of being absorbed by
    # its end-of-paragraph mark.
    $target = $d.Range($paraStart, $paraEnd - 1)
    $target.Text = $newText

    # Re-apply the language formatting the instrText run carried so the
    # resulting run keeps its <w:rPr><w:lang .../></w:rPr>. Use the full
    # paragraph range (including its end-of-paragraph mark) - a range
    # trimmed to just the new text does not reliably stick the
    # character-level language property onto the run.
    for ($p = 1; $p -le $d.Paragraphs.Count; $p++) {
        $pr = $d.Paragraphs.Item($p).Range
        if ($pr.Start -eq $paraStart) {
            $pr.LanguageID = "en-US"
            break
        }
    }
}
